# Windows Generic Credential Functionality Added
# Update the URL and Alphabets values on the Settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# B3: URL value used for login
$ws.Range("B3").Value = "https://masmdvapp.eclinicalweb.com/mobiledoc/jsp/webemr/login/newLogin.jsp"

# B6: Alphabets value (comma separated)
$ws.Range("B6").Value = "ab,ac"

$ws.Range("B18").Select()
